$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.474805333333333
$ws.Range("H2").Value = 4.424416
$ws.Range("I2").Value = 0.0172046112235441
$ws.Range("J2").Value = 0.0172046112235441
$ws.Range("M2").Value = 28.31444233333334
$ws.Range("N2").Value = 84.94332700000001
$ws.Range("O2").Value = 0.2747173016130739
$ws.Range("P2").Value = 0.2747173016130739
$ws.Range("Q2").Value = 41.75829056355911
$ws.Range("R2").Value = 375.824615072032
$ws.Range("S2").Value = 0.004726404370634043
$ws.Range("T2").Value = 0.004726404370634041
$ws.Range("G3").Value = 1.474805333333333
$ws.Range("H3").Value = 4.424416
$ws.Range("I3").Value = 0.0172046112235441
$ws.Range("J3").Value = 0.0172046112235441
$ws.Range("O3").Value = 0.2090339131726295
$ws.Range("P3").Value = 0.2090339131726295
$ws.Range("Q3").Value = 31.77411408981689
$ws.Range("R3").Value = 285.967026808352
$ws.Range("S3").Value = 0.003596347208671165
$ws.Range("T3").Value = 0.003596347208671164
$ws.Range("G4").Value = 1.474805333333333
$ws.Range("H4").Value = 4.424416
$ws.Range("I4").Value = 0.0172046112235441
$ws.Range("J4").Value = 0.0172046112235441
$ws.Range("M4").Value = 5.413469333333334
$ws.Range("N4").Value = 16.240408
$ws.Range("O4").Value = 0.0525235026743817
$ws.Range("P4").Value = 0.0525235026743817
$ws.Range("Q4").Value = 7.983813444636445
$ws.Range("R4").Value = 71.85432100172801
$ws.Range("S4").Value = 0.000903646443611516
$ws.Range("T4").Value = 0.0009036464436115159
$ws.Range("G5").Value = 1.474805333333333
$ws.Range("H5").Value = 4.424416
$ws.Range("I5").Value = 0.0172046112235441
$ws.Range("J5").Value = 0.0172046112235441
$ws.Range("M5").Value = 47.79503400000001
$ws.Range("N5").Value = 143.385102
$ws.Range("O5").Value = 0.4637252825399149
$ws.Range("P5").Value = 0.4637252825399149
$ws.Range("Q5").Value = 70.48837105004802
$ws.Range("R5").Value = 634.3953394504321
$ws.Range("S5").Value = 0.007978213200627379
$ws.Range("T5").Value = 0.007978213200627377
$ws.Range("I6").Value = 0.8384471733397276
$ws.Range("J6").Value = 0.8384471733397275
$ws.Range("M6").Value = 28.31444233333334
$ws.Range("N6").Value = 84.94332700000001
$ws.Range("O6").Value = 0.2747173016130739
$ws.Range("P6").Value = 0.2747173016130739
$ws.Range("Q6").Value = 2035.04282843671
$ws.Range("R6").Value = 18315.38545593039
$ws.Range("S6").Value = 0.2303359450049992
$ws.Range("T6").Value = 0.2303359450049992
$ws.Range("I7").Value = 0.8384471733397276
$ws.Range("J7").Value = 0.8384471733397275
$ws.Range("O7").Value = 0.2090339131726295
$ws.Range("P7").Value = 0.2090339131726295
$ws.Range("S7").Value = 0.1752638936317333
$ws.Range("T7").Value = 0.1752638936317333
$ws.Range("I8").Value = 0.8384471733397276
$ws.Range("J8").Value = 0.8384471733397275
$ws.Range("M8").Value = 5.413469333333334
$ws.Range("N8").Value = 16.240408
$ws.Range("O8").Value = 0.0525235026743817
$ws.Range("P8").Value = 0.0525235026743817
$ws.Range("Q8").Value = 389.0820738783422
$ws.Range("R8").Value = 3501.73866490508
$ws.Range("S8").Value = 0.04403818235123697
$ws.Range("T8").Value = 0.04403818235123696
$ws.Range("I9").Value = 0.8384471733397276
$ws.Range("J9").Value = 0.8384471733397275
$ws.Range("M9").Value = 47.79503400000001
$ws.Range("N9").Value = 143.385102
$ws.Range("O9").Value = 0.4637252825399149
$ws.Range("P9").Value = 0.4637252825399149
$ws.Range("Q9").Value = 3435.17064653903
$ws.Range("R9").Value = 30916.53581885127
$ws.Range("S9").Value = 0.3888091523517582
$ws.Range("T9").Value = 0.3888091523517581
$ws.Range("G10").Value = 12.33901866666667
$ws.Range("H10").Value = 37.017056
$ws.Range("I10").Value = 0.143943077938458
$ws.Range("J10").Value = 0.143943077938458
$ws.Range("M10").Value = 28.31444233333334
$ws.Range("N10").Value = 84.94332700000001
$ws.Range("O10").Value = 0.2747173016130739
$ws.Range("P10").Value = 0.2747173016130739
$ws.Range("Q10").Value = 349.3724324872569
$ws.Range("R10").Value = 3144.351892385312
$ws.Range("S10").Value = 0.03954365395713357
$ws.Range("T10").Value = 0.03954365395713356
$ws.Range("G11").Value = 12.33901866666667
$ws.Range("H11").Value = 37.017056
$ws.Range("I11").Value = 0.143943077938458
$ws.Range("J11").Value = 0.143943077938458
$ws.Range("O11").Value = 0.2090339131726295
$ws.Range("P11").Value = 0.2090339131726295
$ws.Range("Q11").Value = 265.8394148771591
$ws.Range("R11").Value = 2392.554733894432
$ws.Range("S11").Value = 0.03008898485558867
$ws.Range("T11").Value = 0.03008898485558866
$ws.Range("G12").Value = 12.33901866666667
$ws.Range("H12").Value = 37.017056
$ws.Range("I12").Value = 0.143943077938458
$ws.Range("J12").Value = 0.143943077938458
$ws.Range("M12").Value = 5.413469333333334
$ws.Range("N12").Value = 16.240408
$ws.Range("O12").Value = 0.0525235026743817
$ws.Range("P12").Value = 0.0525235026743817
$ws.Range("Q12").Value = 66.79689915542757
$ws.Range("R12").Value = 601.172092398848
$ws.Range("S12").Value = 0.007560394639059331
$ws.Range("T12").Value = 0.00756039463905933
$ws.Range("G13").Value = 12.33901866666667
$ws.Range("H13").Value = 37.017056
$ws.Range("I13").Value = 0.143943077938458
$ws.Range("J13").Value = 0.143943077938458
$ws.Range("M13").Value = 47.79503400000001
$ws.Range("N13").Value = 143.385102
$ws.Range("O13").Value = 0.4637252825399149
$ws.Range("P13").Value = 0.4637252825399149
$ws.Range("Q13").Value = 589.7438166999681
$ws.Range("R13").Value = 5307.694350299712
$ws.Range("S13").Value = 0.06675004448667642
$ws.Range("T13").Value = 0.0667500444866764
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.034729
$ws.Range("H14").Value = 0.104187
$ws.Range("I14").Value = 0.0004051374982703682
$ws.Range("J14").Value = 0.0004051374982703681
$ws.Range("M14").Value = 28.31444233333334
$ws.Range("N14").Value = 84.94332700000001
$ws.Range("O14").Value = 0.2747173016130739
$ws.Range("P14").Value = 0.2747173016130739
$ws.Range("Q14").Value = 0.9833322677943335
$ws.Range("R14").Value = 8.849990410149001
$ws.Range("S14").Value = 0.000111298280307107
$ws.Range("T14").Value = 0.0001112982803071069
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.034729
$ws.Range("H15").Value = 0.104187
$ws.Range("I15").Value = 0.0004051374982703682
$ws.Range("J15").Value = 0.0004051374982703681
$ws.Range("O15").Value = 0.2090339131726295
$ws.Range("P15").Value = 0.2090339131726295
$ws.Range("Q15").Value = 0.7482229574876668
$ws.Range("R15").Value = 6.734006617389
$ws.Range("S15").Value = 0.00008468747663642449
$ws.Range("T15").Value = 0.00008468747663642447
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.034729
$ws.Range("H16").Value = 0.104187
$ws.Range("I16").Value = 0.0004051374982703682
$ws.Range("J16").Value = 0.0004051374982703681
$ws.Range("M16").Value = 5.413469333333334
$ws.Range("N16").Value = 16.240408
$ws.Range("O16").Value = 0.0525235026743817
$ws.Range("P16").Value = 0.0525235026743817
$ws.Range("Q16").Value = 0.1880043764773334
$ws.Range("R16").Value = 1.692039388296
$ws.Range("S16").Value = 0.00002127924047389599
$ws.Range("T16").Value = 0.00002127924047389599
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.034729
$ws.Range("H17").Value = 0.104187
$ws.Range("I17").Value = 0.0004051374982703682
$ws.Range("J17").Value = 0.0004051374982703681
$ws.Range("M17").Value = 47.79503400000001
$ws.Range("N17").Value = 143.385102
$ws.Range("O17").Value = 0.4637252825399149
$ws.Range("P17").Value = 0.4637252825399149
$ws.Range("Q17").Value = 1.659873735786001
$ws.Range("R17").Value = 14.938863622074
$ws.Range("S17").Value = 0.0001878725008529408
$ws.Range("T17").Value = 0.0001878725008529407
